$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume values, preserving original literal-text
# formatting (trailing zeros, padding spaces) the way the source inline-string
# cells stored them. A leading apostrophe forces Excel to keep the assignment
# as literal text instead of re-parsing it as a number (which would silently
# drop significant trailing zeros, e.g. "1.0000" -> 1). Resetting the cell
# Style to "Normal" afterwards removes the quote-prefix/text-format styling
# Excel applies for that trick, so the underlying cell format is untouched.

$ws.Range("D2").Value = "'26.494.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.15%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.733.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.39%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'246.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.39%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E7").Value = "'  +1.85%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2665"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.71%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -0.47%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.736.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.20%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07068"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.92%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'15.73"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'4.614"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.63%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.6104"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.41%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'77.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.17%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.0000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.07%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.485.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.22%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.9999"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.08%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +4.62%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.98%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'1.957.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.35%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.528"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.48%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'8.741"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.49%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'5.250"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.03%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'139.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.47%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'15.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.68%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.781"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.56%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.410"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.06%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'108.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.09%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'3.980"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'0.08050"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.42%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.688"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.32%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.04576"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.50%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.9997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.07%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D36").Value = "'1.013"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.35%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.6368"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.11%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.9058"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.50%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.039"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.40%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.400"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.27%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.002"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.32%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.01510"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.15%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -10.07%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'5.444"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -5.28%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3891"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.14%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'6.963"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +2.54%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -1.93%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05389"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.99%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'30.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.58%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'7.807"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.09%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'1.253"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.41%  "
$ws.Range("E51").Style = "Normal"
